$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1080459066698215
$ws.Range("H2").Value = 62.92134219080646
$ws.Range("G3").Value = 0.07945207144514034
$ws.Range("H3").Value = -32.81740712050722
$ws.Range("G4").Value = -0.3128084454615781
$ws.Range("H4").Value = -13.75044866941243
$ws.Range("G5").Value = -0.3855258351895
$ws.Range("H5").Value = 3.38058300331389
$ws.Range("G6").Value = 0.2192504910394083
$ws.Range("H6").Value = 11.21044951032592
$ws.Range("G7").Value = 0.2632554193539478
$ws.Range("H7").Value = 26.94209929759376
$ws.Range("G8").Value = 0.121587203999592
$ws.Range("H8").Value = 19.32478534167178
$ws.Range("G9").Value = 0.1422309851853292
$ws.Range("H9").Value = 12.45328614854265
$ws.Range("G10").Value = 0.043524765108687
$ws.Range("H10").Value = -29.15637969032777
$ws.Range("G11").Value = 0.01746674402821073
$ws.Range("H11").Value = -65.01764235633378
$ws.Range("G12").Value = 0.09706595934071878
$ws.Range("H12").Value = 4.859980135159919
$ws.Range("G13").Value = 0.09434475753397736
$ws.Range("H13").Value = 23.79947752065303
$ws.Range("G14").Value = 0.2572524436251866
$ws.Range("H14").Value = 13.84783722170463
$ws.Range("G15").Value = 0.2393226478642493
$ws.Range("H15").Value = -2.857543353036252
$ws.Range("G16").Value = 0.152131098523302
$ws.Range("H16").Value = 33.74767110360447
$ws.Range("G17").Value = 0.1335680796103169
$ws.Range("H17").Value = -10.60201898420296
$ws.Range("G18").Value = -0.03317052758298038
$ws.Range("H18").Value = -270.5397225300595
$ws.Range("G19").Value = 0.02587712002003402
$ws.Range("H19").Value = 6.834056682538278
$ws.Range("G20").Value = 0.1190659008633763
$ws.Range("H20").Value = 39.97681390328445
$ws.Range("G21").Value = 0.08353812266868749
$ws.Range("H21").Value = 27.6288168921176
$ws.Range("G22").Value = 0.2196208578153532
$ws.Range("H22").Value = 14.64759510532156
$ws.Range("G23").Value = 0.2216983275189489
$ws.Range("H23").Value = 2.777376565304057
$ws.Range("G24").Value = -0.02485831938478298
$ws.Range("H24").Value = -553.209067490769
$ws.Range("G25").Value = -0.01897896224802626
$ws.Range("H25").Value = 18.39575133120133
$ws.Range("G26").Value = 0.2186390601354116
$ws.Range("H26").Value = 6.722793034212709
$ws.Range("G27").Value = 0.2361584905972731
$ws.Range("H27").Value = 22.43515734709196
$ws.Range("G28").Value = 0.00460185823063228
$ws.Range("H28").Value = -93.12269373827091
$ws.Range("G29").Value = 0.1061640246092083
$ws.Range("H29").Value = 12.62345554563519
$ws.Range("I2").Value = -25.23369985902986
